$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report date (D5): 2026.02.11 06:10 -> 2026.02.13 06:10
$ws.Range("D5").Value = "2026.02.13 06:10"

# New trade history rows appended after row 248 (rows 249-260)
$rows = @(
    @{ r=249; A="2026.02.11 15:22:01"; B=809278101; C="NAS100"; D="buy";  E=10;                 F=25194.5;              I="2026.02.11 16:30:01"; J=25238.19;            K=0; L=0; M=436.9;    N="[tp 25238.19]" },
    @{ r=250; A="2026.02.11 11:05:54"; B=808191793; C="AUDUSD"; D="buy";  E=0.85;               F=0.71094;              I="2026.02.11 16:30:02"; J=0.70896;             K=0; L=0; M=-168.3;   N="[sl 0.70896]" },
    @{ r=251; A="2026.02.11 11:40:20"; B=808309734; C="AUDUSD"; D="buy";  E=0.82;               F=0.7105399999999999;   I="2026.02.11 16:30:08"; J=0.70849;             K=0; L=0; M=-168.1;   N="[sl 0.70849]" },
    @{ r=252; A="2026.02.11 15:11:55"; B=809225290; C="NAS100"; D="buy";  E=8.699999999999999;  F=25194.25;             I="2026.02.11 16:46:12"; J=25307.24;            K=0; L=0; M=983.01;   N="[tp 25307.24]" },
    @{ r=253; A="2026.02.11 11:30:28"; B=808265079; C="LTCUSD"; D="sell"; E=1;                  F=51.38;                I="2026.02.11 17:02:52"; J=53.28;               K=0; L=0; M=-190;     N="[sl 53.28]" },
    @{ r=254; A="2026.02.11 18:35:58"; B=810945618; C="NAS100"; D="sell"; E=8.199999999999999;  F=25081.32;             I="2026.02.11 23:30:40"; J=25208;               K=0; L=0; M=-1038.78; N="close_before_mar" },
    @{ r=255; A="2026.02.12 09:18:57"; B=814015047; C="BTCUSD"; D="sell"; E=0.49;               F=67135.66;             I="2026.02.12 09:19:28"; J=67135.16;            K=0; L=0; M=0.24;     N="[tp 67135.16]" },
    @{ r=256; A="2026.02.12 09:35:59"; B=814090386; C="BTCUSD"; D="sell"; E=0.48;               F=67178.98;             I="2026.02.12 09:36:00"; J=67184.03999999999;   K=0; L=0; M=-2.43;    N="[tp 67184.04]" },
    @{ r=257; A="2026.02.12 11:01:13"; B=814382898; C="XAUUSD"; D="buy";  E=0.08;               F=5064.66;              I="2026.02.12 17:04:05"; J=5071.43;             K=0; L=0; M=54.16;    N="[sl 5071.43]" },
    @{ r=258; A="2026.02.12 15:35:12"; B=815389911; C="NAS100"; D="buy";  E=10;                 F=25303.75;             I="2026.02.12 17:44:10"; J=25250.92;            K=0; L=0; M=-528.3;   N="[sl 25250.92]" },
    @{ r=259; A="2026.02.12 11:08:36"; B=814416556; C="AUDUSD"; D="buy";  E=0.8100000000000001; F=0.71223;              I="2026.02.12 19:15:43"; J=0.71008;             K=0; L=0; M=-174.15;  N="[sl 0.71008]" },
    @{ r=260; A="2026.02.12 11:00:13"; B=814378525; C="AUDUSD"; D="buy";  E=0.76;               F=0.71216;              I="2026.02.12 19:16:17"; J=0.70987;             K=0; L=0; M=-174.04;  N="[sl 0.70987]" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
}
